$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 03:14"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8453850
$ws.Range("C4").Value = 54524
$ws.Range("D4").Value = 5496779
$ws.Range("E4").Value = 2731859
$ws.Range("G4").Value = 432
$ws.Range("H4").Value = 225212

# Row 86 - Australia
$ws.Range("D86").Value = 25112
$ws.Range("E86").Value = 1382

# Row 138 - Republica de Africa Central
$ws.Range("B138").Value = 4856
$ws.Range("C138").Value = 1
$ws.Range("E138").Value = 2870

# Row 166 - Liberia
$ws.Range("B166").Value = 1381
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 1271
$ws.Range("E166").Value = 28

# Row 171 - San Martin (Parte Holandesa)
$ws.Range("B171").Value = 756
$ws.Range("C171").Value = 3
$ws.Range("D171").Value = 675
$ws.Range("E171").Value = 59

# Row 174 - Islas Turcas y Caicos
$ws.Range("D174").Value = 689
$ws.Range("E174").Value = 3

# Row 185 - Isla de Man
$ws.Range("D185").Value = 321
$ws.Range("E185").Value = 3

# Row 190 - Islas Caimanes
$ws.Range("B190").Value = 235
$ws.Range("C190").Value = 2
$ws.Range("E190").Value = 22
